$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the Task 7 placeholder with the full task description.
$ws.Range("A20").Value2 = "Task 7 Description: Calculate the multiplication of arccos(0.5) and arccos(-0.5) using the acos(x) Function ((acos(0.5)*acos(-0.5))"

# Update the current selection/view to match the saved workbook state.
$ws.Range("B22").Select()
